$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'89.169.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.36%  "

# Row 3
$ws.Range("D3").Value = "'3.096.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.61%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "'212.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.86%  "

# Row 6
$ws.Range("D6").Value = "'622.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.17%  "

# Row 7
$ws.Range("D7").Value = "'0.374"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -7.59%  "

# Row 8
$ws.Range("D8").Value = "'0.819"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +14.94%  "

# Row 9
$ws.Range("D9").Value = "'0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.04%  "

# Row 10
$ws.Range("D10").Value = "'3.096.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.67%  "

# Row 11
$ws.Range("D11").Value = "'0.624"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.51%  "

# Row 12
$ws.Range("E12").Value = "  -0.91%  "

# Row 13
$ws.Range("E13").Value = "  -7.87%  "

# Row 14
$ws.Range("D14").Value = "'5.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.72%  "

# Row 15
$ws.Range("D15").Value = "'88.857.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.24%  "

# Row 16
$ws.Range("D16").Value = "'32.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.07%  "

# Row 17
$ws.Range("D17").Value = "'3.664.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.61%  "

# Row 18
$ws.Range("D18").Value = "'3.095.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.92%  "

# Row 19
$ws.Range("D19").Value = "'3.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.51%  "

# Row 20
$ws.Range("D20").Value = "'0.0000212"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.27%  "

# Row 21
$ws.Range("D21").Value = "'13.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.62%  "

# Row 22
$ws.Range("D22").Value = "'423.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.50%  "

# Row 23
$ws.Range("D23").Value = "'8.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.48%  "

# Row 24
$ws.Range("D24").Value = "'4.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.93%  "

# Row 25
$ws.Range("D25").Value = "'5.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.13%  "

# Row 26
$ws.Range("D26").Value = "'11.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.74%  "

# Row 27
$ws.Range("D27").Value = "'82.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.70%  "

# Row 28
$ws.Range("D28").Value = "'3.235.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.06%  "

# Row 29
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.07%  "

# Row 30
$ws.Range("D30").Value = "'0.171"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.20%  "

# Row 31
$ws.Range("D31").Value = "'1.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.35%  "

# Row 32
$ws.Range("D32").Value = "'8.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.46%  "

# Row 33
$ws.Range("D33").Value = "'510.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.41%  "

# Row 34
$ws.Range("E34").Value = "  -13.53%  "

# Row 35
$ws.Range("E35").Value = "  -5.44%  "

# Row 36
$ws.Range("E36").Value = "  -3.99%  "

# Row 37
$ws.Range("E37").Value = "  -5.95%  "

# Row 38
$ws.Range("D38").Value = "'22.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.47%  "

# Row 39
$ws.Range("D39").Value = "'0.131"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.01%  "

# Row 40
$ws.Range("D40").Value = "'22.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.59%  "

# Row 41
$ws.Range("E41").Value = "  +0.27%  "

# Row 42
$ws.Range("E42").Value = "  -0.02%  "

# Row 43
$ws.Range("E43").Value = "  -3.34%  "

# Row 44
$ws.Range("E44").Value = "  -7.20%  "

# Row 45
$ws.Range("D45").Value = "'145.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.32%  "

# Row 46
$ws.Range("E46").Value = "  +4.09%  "

# Row 47
$ws.Range("D47").Value = "'0.0695"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +12.86%  "

# Row 48
$ws.Range("D48").Value = "'43.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.08%  "

# Row 49
$ws.Range("D49").Value = "'161.13"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.42%  "

# Row 50
$ws.Range("D50").Value = "'1.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.28%  "

# Row 51
$ws.Range("E51").Value = "  -6.62%  "

